$d = $word.ActiveDocument

function Add-StyledParagraph($text, $styleName) {
    $endPos = $d.Content.End
    $r = $d.Range($endPos, $endPos)
    $r.InsertAfter("`r" + $text)
    $p = $d.Paragraphs.Last
    $p.Range.Font.Reset()
    $p.Style = $styleName
}

# Remember where the new glossary content starts (right after the existing
# bibliography, i.e. at the current end of the document body).
$glossaryStart = $d.Content.End

Add-StyledParagraph "Glossary" "Ttulo1"
Add-StyledParagraph "Some term" "DefinitionTerm"
Add-StyledParagraph "The definition of this term" "Definition"
Add-StyledParagraph "Another term" "DefinitionTerm"
Add-StyledParagraph "The definition of this other term" "Definition"
Add-StyledParagraph "Get it?" "DefinitionTerm"
Add-StyledParagraph "Another term with its corresponding definitions" "Definition"

$glossaryEnd = $d.Content.End
$glossaryRange = $d.Range($glossaryStart, $glossaryEnd)
$d.Bookmarks.Add("glossary", $glossaryRange)

Write-Output "Glossary added."
